# cleaned up and merged login tests for both partner and insider apps
#
# Adds a new "setUp" worksheet (holding the partner/insider QA URLs used by
# the merged login tests) as the last sheet in the workbook, with a
# hyperlink on the insider URL, and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$setupSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$setupSheet.Name = "setUp"

# Values: header + the two QA environment URLs.
$setupSheet.Range("A1").Value = "url"
$setupSheet.Range("A2").Value = "https://partner.qa.upnorway.net/"
$setupSheet.Range("A3").Value = "https://insider.qa.upnorway.net/"

# Hyperlink the insider URL cell and give it the built-in Hyperlink style.
$setupSheet.Hyperlinks.Add($setupSheet.Range("A3"), "https://insider.qa.upnorway.net/")
$setupSheet.Range("A3").Style = "Hyperlink"

# Make the new sheet the active tab / selection, matching the saved view state.
$setupSheet.Activate()
$setupSheet.Range("L11").Select()
